$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 626 (shifts existing rows 626..667 down to 627..668)
$ws.Rows.Item(626).Insert()

# Populate the newly inserted row 626 with the new data point.
# Column A holds a date-formatted-looking string that the sheet stores as
# plain text (matching every other row), so force text interpretation
# while assigning it to avoid Excel auto-converting it into a real date
# serial number, then restore the default "Normal" style so the cell ends
# up unstyled like its neighbours.
$ws.Cells.Item(626, 1).NumberFormat = "@"
$ws.Cells.Item(626, 1).Value = "2026/01/11"
$ws.Cells.Item(626, 1).Style = "Normal"
$ws.Cells.Item(626, 2).Value = "日"
$ws.Cells.Item(626, 3).Value = 8
$ws.Cells.Item(626, 4).Value = 151
